$d = $word.ActiveDocument
$c = $d.Content
$c.SetRange($c.End, $c.End)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:t>Get-Process node -</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ErrorAction</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SilentlyContinue</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> | Stop-Process -Force -</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ErrorAction</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SilentlyContinue</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>Start-Sleep -Seconds 3</w:t></w:r><w:r><w:t>0</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve"># Clear </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> cache for client</w:t></w:r></w:p><w:p><w:r><w:t>cd "C:\Users\stuat\Documents\psychic-chat-poc\client"</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Remove-Item -Path </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>node_modules</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramStart"/><w:r><w:t>\.cache</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> -Recurse -Force -</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ErrorAction</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SilentlyContinue</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>Write-Output "Cleared React cache"</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$c.InsertXML($xml)
